# Update the "complete_Table2" worksheet so that the error-covariance /
# squared-error vectors (columns L and M) for the RDCC rows are filled in,
# and the K-column figures that depend on them are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- K column: small refinements to the CP-copula figures -----------------
$ws.Range("K8").Value = 0.014
$ws.Range("K9").Value = 0.986

# --- New L/M (vector a / vector b) entries in the parameter block ---------
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0.005

$ws.Range("L11").Value = 0.02
$ws.Range("M11").Value = 0.005

$ws.Range("L12").Value = 0.52

$ws.Range("L13").Value = 0.635

$ws.Range("M14").Value = 0.495

# --- LL Decomposition rows: refreshed Total_LL for RDCC_CP and new lower
#     bounds (L/M) used as penalty values for the squared errors / passed
#     conditional covariances -----------------------------------------------
foreach ($r in 17..20) {
    $ws.Cells.Item($r, 11).Value = -8821.743   # column K
    $ws.Cells.Item($r, 12).Value = -10000000   # column L
    $ws.Cells.Item($r, 13).Value = -10000000   # column M
}

# --- Column M got a touch wider to fit the new figures --------------------
$ws.Columns.Item(13).ColumnWidth = 8.8
